$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

function Set-PlainCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

$ws1 = $wb.Worksheets.Item("Sheet1")
Set-PlainCell $ws1 7 1 "Andrew Armstrong"
Set-TextCell  $ws1 7 2 "09/23/2023"
Set-TextCell  $ws1 7 3 "6"
Set-TextCell  $ws1 7 4 "84"
Set-PlainCell $ws1 7 5 "CH"
Set-PlainCell $ws1 7 6 "Ball"

Set-PlainCell $ws1 8 1 "Andrew Armstrong"
Set-TextCell  $ws1 8 2 "09/23/2023"
Set-TextCell  $ws1 8 3 "7"
Set-TextCell  $ws1 8 4 "88"
Set-PlainCell $ws1 8 5 "SL"
Set-PlainCell $ws1 8 6 "HBP"

Set-PlainCell $ws1 9 1 "Andrew Armstrong"
Set-TextCell  $ws1 9 2 "09/23/2023"
Set-TextCell  $ws1 9 3 "8"
Set-TextCell  $ws1 9 4 "90"
Set-PlainCell $ws1 9 5 "FB"
Set-PlainCell $ws1 9 6 "Strike looking"

Set-PlainCell $ws1 10 1 "Andrew Armstrong"
Set-TextCell  $ws1 10 2 "09/23/2023"
Set-TextCell  $ws1 10 3 "9"
Set-TextCell  $ws1 10 4 "92"
Set-PlainCell $ws1 10 5 "FB"
Set-PlainCell $ws1 10 6 "Foul Ball"

Set-PlainCell $ws1 11 1 "Andrew Armstrong"
Set-TextCell  $ws1 11 2 "09/23/2023"
Set-TextCell  $ws1 11 3 "10"
Set-TextCell  $ws1 11 4 "90"
Set-PlainCell $ws1 11 5 "FB"
Set-PlainCell $ws1 11 6 "Ball"

Set-PlainCell $ws1 12 1 "Andrew Armstrong"
Set-TextCell  $ws1 12 2 "09/23/2023"
Set-TextCell  $ws1 12 3 "11"
Set-TextCell  $ws1 12 4 "77"
Set-PlainCell $ws1 12 5 "CB"
Set-PlainCell $ws1 12 6 "Ball"

Set-PlainCell $ws1 13 1 "Andrew Armstrong"
Set-TextCell  $ws1 13 2 "09/23/2023"
Set-TextCell  $ws1 13 3 "12"
Set-TextCell  $ws1 13 4 "88"
Set-PlainCell $ws1 13 5 "FB"
Set-PlainCell $ws1 13 6 "Strikeout swinging"

Set-PlainCell $ws1 14 1 "Andrew Armstrong"
Set-TextCell  $ws1 14 2 "09/23/2023"
Set-TextCell  $ws1 14 3 "13"
Set-TextCell  $ws1 14 4 "90"
Set-PlainCell $ws1 14 5 "FB"
Set-PlainCell $ws1 14 6 "Hit"

Set-PlainCell $ws1 15 1 "Andrew Armstrong"
Set-TextCell  $ws1 15 2 "09/23/2023"
Set-TextCell  $ws1 15 3 "14"
Set-TextCell  $ws1 15 4 "88"
Set-PlainCell $ws1 15 5 "Knuck"
Set-PlainCell $ws1 15 6 "Ball"

Set-PlainCell $ws1 16 1 "Andrew Armstrong"
Set-TextCell  $ws1 16 2 "09/23/2023"
Set-TextCell  $ws1 16 3 "15"
Set-TextCell  $ws1 16 4 "90"
Set-PlainCell $ws1 16 5 "FB"
Set-PlainCell $ws1 16 6 "Ball"

Set-PlainCell $ws1 17 1 "Andrew Armstrong"
Set-TextCell  $ws1 17 2 "09/23/2023"
Set-TextCell  $ws1 17 3 "16"
Set-TextCell  $ws1 17 4 "88"
Set-PlainCell $ws1 17 5 "FB"
Set-PlainCell $ws1 17 6 "Ball"

Set-PlainCell $ws1 18 1 "Andrew Armstrong"
Set-TextCell  $ws1 18 2 "09/23/2023"
Set-TextCell  $ws1 18 3 "17"
Set-TextCell  $ws1 18 4 "75"
Set-PlainCell $ws1 18 5 "Knuck"
Set-PlainCell $ws1 18 6 "Walk"

$ws2 = $wb.Worksheets.Item("pitch breakdown")
Set-PlainCell $ws2 7 1 "Andrew Armstrong"
Set-TextCell  $ws2 7 2 "09/23/2023"
Set-TextCell  $ws2 7 3 "6"
Set-TextCell  $ws2 7 4 "84"
Set-PlainCell $ws2 7 5 "CH"
Set-PlainCell $ws2 7 6 "Ball"
Set-PlainCell $ws2 7 7 "Ball"
Set-PlainCell $ws2 7 8 "No swing"
Set-PlainCell $ws2 7 9 "nothing"

Set-PlainCell $ws2 8 1 "Andrew Armstrong"
Set-TextCell  $ws2 8 2 "09/23/2023"
Set-TextCell  $ws2 8 3 "7"
Set-TextCell  $ws2 8 4 "88"
Set-PlainCell $ws2 8 5 "SL"
Set-PlainCell $ws2 8 6 "HBP"
Set-PlainCell $ws2 8 7 "Ball"
Set-PlainCell $ws2 8 8 "No swing"
Set-PlainCell $ws2 8 9 "free base"

Set-PlainCell $ws2 9 1 "Andrew Armstrong"
Set-TextCell  $ws2 9 2 "09/23/2023"
Set-TextCell  $ws2 9 3 "8"
Set-TextCell  $ws2 9 4 "90"
Set-PlainCell $ws2 9 5 "FB"
Set-PlainCell $ws2 9 6 "Strike looking"
Set-PlainCell $ws2 9 7 "Strike"
Set-PlainCell $ws2 9 8 "No swing"
Set-PlainCell $ws2 9 9 "nothing"

Set-PlainCell $ws2 10 1 "Andrew Armstrong"
Set-TextCell  $ws2 10 2 "09/23/2023"
Set-TextCell  $ws2 10 3 "9"
Set-TextCell  $ws2 10 4 "92"
Set-PlainCell $ws2 10 5 "FB"
Set-PlainCell $ws2 10 6 "Foul Ball"
Set-PlainCell $ws2 10 7 "Strike"
Set-PlainCell $ws2 10 8 "Swing contact"
Set-PlainCell $ws2 10 9 "nothing"

Set-PlainCell $ws2 11 1 "Andrew Armstrong"
Set-TextCell  $ws2 11 2 "09/23/2023"
Set-TextCell  $ws2 11 3 "10"
Set-TextCell  $ws2 11 4 "90"
Set-PlainCell $ws2 11 5 "FB"
Set-PlainCell $ws2 11 6 "Ball"
Set-PlainCell $ws2 11 7 "Ball"
Set-PlainCell $ws2 11 8 "No swing"
Set-PlainCell $ws2 11 9 "nothing"

Set-PlainCell $ws2 12 1 "Andrew Armstrong"
Set-TextCell  $ws2 12 2 "09/23/2023"
Set-TextCell  $ws2 12 3 "11"
Set-TextCell  $ws2 12 4 "77"
Set-PlainCell $ws2 12 5 "CB"
Set-PlainCell $ws2 12 6 "Ball"
Set-PlainCell $ws2 12 7 "Ball"
Set-PlainCell $ws2 12 8 "No swing"
Set-PlainCell $ws2 12 9 "nothing"

Set-PlainCell $ws2 13 1 "Andrew Armstrong"
Set-TextCell  $ws2 13 2 "09/23/2023"
Set-TextCell  $ws2 13 3 "12"
Set-TextCell  $ws2 13 4 "88"
Set-PlainCell $ws2 13 5 "FB"
Set-PlainCell $ws2 13 6 "Strikeout swinging"
Set-PlainCell $ws2 13 7 "Strike"
Set-PlainCell $ws2 13 8 "Swing no contact"
Set-PlainCell $ws2 13 9 "not free base"

Set-PlainCell $ws2 14 1 "Andrew Armstrong"
Set-TextCell  $ws2 14 2 "09/23/2023"
Set-TextCell  $ws2 14 3 "13"
Set-TextCell  $ws2 14 4 "90"
Set-PlainCell $ws2 14 5 "FB"
Set-PlainCell $ws2 14 6 "Hit"
Set-PlainCell $ws2 14 7 "Strike"
Set-PlainCell $ws2 14 8 "Swing contact"
Set-PlainCell $ws2 14 9 "not free base"

Set-PlainCell $ws2 15 1 "Andrew Armstrong"
Set-TextCell  $ws2 15 2 "09/23/2023"
Set-TextCell  $ws2 15 3 "14"
Set-TextCell  $ws2 15 4 "88"
Set-PlainCell $ws2 15 5 "Knuck"
Set-PlainCell $ws2 15 6 "Ball"
Set-PlainCell $ws2 15 7 "Ball"
Set-PlainCell $ws2 15 8 "No swing"
Set-PlainCell $ws2 15 9 "nothing"

Set-PlainCell $ws2 16 1 "Andrew Armstrong"
Set-TextCell  $ws2 16 2 "09/23/2023"
Set-TextCell  $ws2 16 3 "15"
Set-TextCell  $ws2 16 4 "90"
Set-PlainCell $ws2 16 5 "FB"
Set-PlainCell $ws2 16 6 "Ball"
Set-PlainCell $ws2 16 7 "Ball"
Set-PlainCell $ws2 16 8 "No swing"
Set-PlainCell $ws2 16 9 "nothing"

Set-PlainCell $ws2 17 1 "Andrew Armstrong"
Set-TextCell  $ws2 17 2 "09/23/2023"
Set-TextCell  $ws2 17 3 "16"
Set-TextCell  $ws2 17 4 "88"
Set-PlainCell $ws2 17 5 "FB"
Set-PlainCell $ws2 17 6 "Ball"
Set-PlainCell $ws2 17 7 "Ball"
Set-PlainCell $ws2 17 8 "No swing"
Set-PlainCell $ws2 17 9 "nothing"

Set-PlainCell $ws2 18 1 "Andrew Armstrong"
Set-TextCell  $ws2 18 2 "09/23/2023"
Set-TextCell  $ws2 18 3 "17"
Set-TextCell  $ws2 18 4 "75"
Set-PlainCell $ws2 18 5 "Knuck"
Set-PlainCell $ws2 18 6 "Walk"
Set-PlainCell $ws2 18 7 "Ball"
Set-PlainCell $ws2 18 8 "No swing"
Set-PlainCell $ws2 18 9 "free base"

$ws3 = $wb.Worksheets.Item("pitcher breakdown")
$ws3.Cells.Item(2, 3).Value = 92
$ws3.Cells.Item(2, 4).Value = 0.5
$ws3.Cells.Item(2, 6).Value = 0.3
$ws3.Cells.Item(2, 8).Value = 0.3
$ws3.Cells.Item(2, 9).Value = 2
